$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: swap the F:V contents of two rows (A:E - index/pais/torneio/
# temporada/data_partida - stay put, only the match data columns move).
# ---------------------------------------------------------------------------
function Swap-MatchRows($rowA, $rowB) {
    $cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
    foreach ($col in $cols) {
        $valA = $ws.Range($col + $rowA).Value2
        $valB = $ws.Range($col + $rowB).Value2
        $ws.Range($col + $rowA).Value = $valB
        $ws.Range($col + $rowB).Value = $valA
    }
}

Swap-MatchRows 36 37
Swap-MatchRows 40 41
Swap-MatchRows 48 49
Swap-MatchRows 70 71

# ---------------------------------------------------------------------------
# Append three new match rows (78, 79, 80) at the end of the sheet.
# ---------------------------------------------------------------------------
function Add-MatchRow($r, $idx, $dateSerial, $home, $homeGoals, $away, $awayGoals, `
                       $hOpenOdds, $hOpenDt, $hCloseOdds, $hCloseDt, `
                       $dOpenOdds, $dOpenDt, $dCloseOdds, $dCloseDt, `
                       $aOpenOdds, $aOpenDt, $aCloseOdds, $aCloseDt, $url) {

    $ws.Range("A36").Copy()
    $ws.Range("A" + $r).PasteSpecial(-4122)
    $ws.Range("A" + $r).Value = $idx

    $ws.Range("B" + $r).Value = "switzerland"
    $ws.Range("C" + $r).Value = "super-league"
    $ws.Range("D" + $r).Value = "2023-2024"

    $ws.Range("E36").Copy()
    $ws.Range("E" + $r).PasteSpecial(-4122)
    $ws.Range("E" + $r).Value = $dateSerial

    $ws.Range("F" + $r).Value = $home
    $ws.Range("G" + $r).Value = $homeGoals
    $ws.Range("H" + $r).Value = $away
    $ws.Range("I" + $r).Value = $awayGoals

    $ws.Range("J" + $r).Value = $hOpenOdds
    $ws.Range("K" + $r).Value = $hOpenDt
    $ws.Range("L" + $r).Value = $hCloseOdds
    $ws.Range("M" + $r).Value = $hCloseDt

    $ws.Range("N" + $r).Value = $dOpenOdds
    $ws.Range("O" + $r).Value = $dOpenDt
    $ws.Range("P" + $r).Value = $dCloseOdds
    $ws.Range("Q" + $r).Value = $dCloseDt

    $ws.Range("R" + $r).Value = $aOpenOdds
    $ws.Range("S" + $r).Value = $aOpenDt
    $ws.Range("T" + $r).Value = $aCloseOdds
    $ws.Range("U" + $r).Value = $aCloseDt

    $ws.Range("V" + $r).Value = $url
}

Add-MatchRow 78 77 45241.75 "Yverdon" 2 "Lausanne" 2 `
    3.38 "05/11/2023 16:42" 3.74 "11/11/2023 17:59" `
    3.81 "05/11/2023 16:42" 3.87 "11/11/2023 17:59" `
    2.06 "05/11/2023 16:42" 1.96 "11/11/2023 17:53" `
    "https://www.betexplorer.com/football/switzerland/super-league/yverdon-lausanne/hjio1qWf/"

Add-MatchRow 79 78 45241.75 "St. Gallen" 4 "Winterthur" 2 `
    1.47 "05/11/2023 16:42" 1.53 "11/11/2023 17:52" `
    4.98 "05/11/2023 16:42" 4.93 "11/11/2023 17:55" `
    5.94 "05/11/2023 16:42" 5.57 "11/11/2023 17:55" `
    "https://www.betexplorer.com/football/switzerland/super-league/st-gallen-winterthur/0rZw351r/"

Add-MatchRow 80 79 45241.85416666666 "Young Boys" 6 "Luzern" 1 `
    1.6 "05/11/2023 14:42" 1.86 "11/11/2023 20:26" `
    4.55 "05/11/2023 14:42" 4.08 "11/11/2023 20:26" `
    4.93 "05/11/2023 14:42" 3.94 "11/11/2023 20:26" `
    "https://www.betexplorer.com/football/switzerland/super-league/young-boys-luzern/bTis2PGl/"
